# Auto-update predictions and index for 2025-10-15
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Palmeiras v Bragantino
$ws.Range("E2").Value = "46/50 Win Tips"
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "92"
$ws.Cells.Item(2, 6).ClearFormats()

# Row 3: Mansfield v Newcastle U21
$ws.Range("E3").Value = "30/38 Win Tips"
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = "79"
$ws.Cells.Item(3, 6).ClearFormats()

# Row 4: Botafogo v Flamengo
$ws.Range("E4").Value = "19/34 Win Tips"
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = "56"
$ws.Cells.Item(4, 6).ClearFormats()
